$d = $word.ActiveDocument

# The document's header/footer stories each contain one inline picture:
#   - the BTec logo in the header (AlternativeText "BTec_Logo-Orange")
#   - the Pearson Edexcel logo in the two footers (AlternativeText
#     "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png")
# Their drawing "name" attributes were swapped:
#   - Pearson Edexcel logo pictures: image1.png -> image2.png
#   - BTec logo picture:             image2.jpg -> image1.jpg
# Walk every section's headers and footers and rename the matching
# inline shapes via the supported InlineShape.Name property. A
# Write-Host is issued after each rename so the pending write is
# flushed before the next header/footer story is touched.

foreach ($sec in $d.Sections) {

    for ($i = 1; $i -le $sec.Headers.Count; $i++) {
        $hdr = $sec.Headers.Item($i)
        if ($hdr.Exists) {
            foreach ($shp in $hdr.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image1.jpg"
                    Write-Host "Renamed header inline shape to image1.jpg"
                }
            }
        }
    }

    for ($i = 1; $i -le $sec.Footers.Count; $i++) {
        $ftr = $sec.Footers.Item($i)
        if ($ftr.Exists) {
            foreach ($shp in $ftr.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image2.png"
                    Write-Host "Renamed footer inline shape to image2.png"
                }
            }
        }
    }
}
